$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 47: Normex s.r.o. / Roman Toda ---
$ws.Cells.Item(47,1).Value = "NORM"
$ws.Cells.Item(47,2).Value = "Roman"
$ws.Cells.Item(47,3).Value = "Toda"
$ws.Cells.Item(47,4).Value = "toda@digitaldocuments.org"
$ws.Cells.Item(47,5).Value = "Normex s.r.o."
$ws.Cells.Item(47,6).Value = 43502

# --- Row 48: PDF Association / Duff Johnson ---
$ws.Cells.Item(48,1).Value = "pdfa"
$ws.Cells.Item(48,2).Value = "Duff"
$ws.Cells.Item(48,3).Value = "Johnson"
$ws.Cells.Item(48,4).Value = "duff.johnson@pdfa.org"
$ws.Cells.Item(48,5).Value = "PDF Association"
$ws.Cells.Item(48,6).Value = 43502

# --- Hyperlinks on the email column (D) ---
$ws.Hyperlinks.Add($ws.Cells.Item(47,4), "mailto:toda@digitaldocuments.org", "")
$ws.Hyperlinks.Add($ws.Cells.Item(48,4), "mailto:duff.johnson@pdfa.org", "")

# --- Update view: scroll/selection to match post-edit state ---
$ws.Range("F47:F48").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F47:F48").Select()
